$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange
$sub = $tr.Characters(1, 7)
$sub.Text = "Lab 10: "
Write-Host "Full text: [$($tr.Text)] len=$($tr.Length)"
Write-Host "RunCount:" $tr.Runs().Count
for ($i = 1; $i -le $tr.Runs().Count; $i++) {
    $r = $tr.Runs($i)
    Write-Host "$i : [$($r.Text)] bold=$($r.Font.Bold)"
}
